$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting existing data (D:K) to (E:L)
$ws.Columns("D:D").Insert()

# Copy formatting from column E (the shifted original column D) into new column D
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)

# Populate new column D with the new fiscal-year data
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 172100
$ws.Range("D9").Value = 160100
$ws.Range("D10").Value = 12100
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = -21700
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 155000
$ws.Range("D18").Value = 17100
$ws.Range("D20").Value = 500
$ws.Range("D21").Value = 27000
$ws.Range("D22").Value = 600
$ws.Range("D23").Value = 17100
$ws.Range("D24").Value = -3300
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 20300
$ws.Range("D27").Value = 20300
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -500
$ws.Range("D33").Value = 20300
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 20300
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 6700
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 108700
$ws.Range("D44").Value = 39400
$ws.Range("D45").Value = 4800
$ws.Range("D46").Value = 159500
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 103400
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 8400
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 271400
$ws.Range("D57").Value = 19800
$ws.Range("D58").Value = 400
$ws.Range("D59").Value = 11300
$ws.Range("D60").Value = 31500
$ws.Range("D61").Value = 12300
$ws.Range("D62").Value = 9000
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 52800
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 101200
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 218600
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 20300
$ws.Range("D83").Value = 9300
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = -18400
$ws.Range("D91").Value = -3800
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -27900
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 9300
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = -37000
$ws.Range("D12").Value = "NA"
$ws.Range("D49").Value = "NA"
